$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: replace the existing "Admin_Report()" test-case mapping with the
# corrected function name, and add a new IT mapping in column E.
$ws.Range("C13").Value = "Admin_Max Request_Report()"
$ws.Range("E13").Value = "IT_11"

# Rows 14-15: add new IT mapping values in column E (previously empty).
$ws.Range("E14").Value = "IT_12"
$ws.Range("E15").Value = "IT_13"

# These four cells pick up an explicit (Automatic/theme) font color, which is
# how they end up differently formatted from the surrounding "s=3" cells.
$ws.Range("C13").Font.ThemeColor = 1
$ws.Range("E13").Font.ThemeColor = 1
$ws.Range("E14").Font.ThemeColor = 1
$ws.Range("E15").Font.ThemeColor = 1

# Final selection left on D16.
$ws.Range("D16").Select()
